$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit reshuffles species-observation data among rows 2,3,4,5,6,10,11,12
# (Taxonsorteringsordning column B also bumps by +14 on every row).
# Build the new values per target row and write them.

$rows = @{
    2  = @{ A = 112243569; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404751; R = 6707073 }
    3  = @{ A = 112243589; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404762; R = 6707097 }
    4  = @{ A = 112243600; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404725; R = 6707036 }
    5  = @{ A = 112243563; B = 89573; D = "NT"; E = 5442; F = "Tallticka";      G = "Porodaedalea pini";       H = "(Brot.) Murrill";      Q = 404744; R = 6707084 }
    6  = @{ A = 112243594; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404742; R = 6706992 }
    10 = @{ A = 112243565; B = 90814; D = "LC"; E = 4364; F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum";   H = "(Fr.:Fr.) P. Karst.";  Q = 404459; R = 6706753 }
    11 = @{ A = 112243588; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404452; R = 6706739 }
    12 = @{ A = 112243573; B = 77650; D = "NT"; E = 6425; F = "Garnlav";        G = "Alectoria sarmentosa";    H = "(Ach.) Ach.";          Q = 404477; R = 6706766 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value  = $vals.A   # A: Id
    $ws.Cells.Item($r, 2).Value  = $vals.B   # B: Taxonsorteringsordning
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Rödlistade
    $ws.Cells.Item($r, 5).Value  = $vals.E   # E: TaxonId
    $ws.Cells.Item($r, 6).Value  = $vals.F   # F: Artnamn
    $ws.Cells.Item($r, 7).Value  = $vals.G   # G: Vetenskapligt namn
    $ws.Cells.Item($r, 8).Value  = $vals.H   # H: Auktor
    $ws.Cells.Item($r, 17).Value = $vals.Q   # Q: Ost
    $ws.Cells.Item($r, 18).Value = $vals.R   # R: Nord
}
